# CH allow valueString as well as valueQuantity
# Update the "Value Types" column (H) for the VA.MHV.PHR.chTest row (row 2)
# from "Quantityĵ" to "Quantityĵ, stringĵ"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observations")

$ws.Range("H2").Value = "Quantityĵ, stringĵ"
